# Horarios actualizados Linea 141 - 768
# Applies the 04:46:27 scrape refresh to all three schedule sheets:
#   LP1912, LP1912-215, 6203-6173
$wb = $excel.ActiveWorkbook

# ===== Sheet 1: LP1912 =====
$ws = $wb.Worksheets.Item(1)

# Refresh timestamp / row-count banner
$ws.Cells.Item(2, "A").Value = 'Última actualización: 04:46:27'
$ws.Cells.Item(3, "A").Value = 'Total filas: 32'

# Data rows (Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada)
# row 18: 215A_EL PATO
$ws.Cells.Item(18, "A").Value = '04:46:27'
$ws.Cells.Item(18, "B").Value = '04:46'
$ws.Cells.Item(18, "C").Value = '215A_EL PATO'
$ws.Cells.Item(18, "D").Value = 0
$ws.Cells.Item(18, "E").Value = 'LP1912'

# row 19: 11_ETCHEVERRY
$ws.Cells.Item(19, "A").Value = '04:46:27'
$ws.Cells.Item(19, "B").Value = '04:53'
$ws.Cells.Item(19, "C").Value = '11_ETCHEVERRY'
$ws.Cells.Item(19, "D").Value = 7
$ws.Cells.Item(19, "E").Value = 'LP1912'

# row 21: 17_ROMERO
$ws.Cells.Item(21, "A").Value = '04:46:27'
$ws.Cells.Item(21, "B").Value = '05:16'
$ws.Cells.Item(21, "C").Value = '17_ROMERO'
$ws.Cells.Item(21, "D").Value = 30
$ws.Cells.Item(21, "E").Value = 'LP1912'

# row 22: 23_HERNANDEZ
$ws.Cells.Item(22, "A").Value = '04:46:27'
$ws.Cells.Item(22, "B").Value = '05:22'
$ws.Cells.Item(22, "C").Value = '23_HERNANDEZ'
$ws.Cells.Item(22, "D").Value = 36
$ws.Cells.Item(22, "E").Value = 'LP1912'

# row 23: 215B_EL PATO
$ws.Cells.Item(23, "A").Value = '04:46:27'
$ws.Cells.Item(23, "B").Value = '05:34'
$ws.Cells.Item(23, "C").Value = '215B_EL PATO'
$ws.Cells.Item(23, "D").Value = 48
$ws.Cells.Item(23, "E").Value = 'LP1912'

# row 24: 215B_EL PATO
$ws.Cells.Item(24, "A").Value = '04:24:09'
$ws.Cells.Item(24, "B").Value = '05:35'
$ws.Cells.Item(24, "C").Value = '215B_EL PATO'
$ws.Cells.Item(24, "D").Value = 71
$ws.Cells.Item(24, "E").Value = 'LP1912'

# row 25: 14_ABASTO
$ws.Cells.Item(25, "A").Value = '03:57:17'
$ws.Cells.Item(25, "B").Value = '05:39'
$ws.Cells.Item(25, "C").Value = '14_ABASTO'
$ws.Cells.Item(25, "D").Value = 102
$ws.Cells.Item(25, "E").Value = 'LP1912'

# row 26: 15_ABASTO
$ws.Cells.Item(26, "A").Value = '04:46:27'
$ws.Cells.Item(26, "B").Value = '05:46'
$ws.Cells.Item(26, "C").Value = '15_ABASTO'
$ws.Cells.Item(26, "D").Value = 60
$ws.Cells.Item(26, "E").Value = 'LP1912'

# row 27: 215A_EL PATO
$ws.Cells.Item(27, "A").Value = '04:24:09'
$ws.Cells.Item(27, "B").Value = '05:48'
$ws.Cells.Item(27, "C").Value = '215A_EL PATO'
$ws.Cells.Item(27, "D").Value = 84
$ws.Cells.Item(27, "E").Value = 'LP1912'

# row 28: 10_OLMOS
$ws.Cells.Item(28, "A").Value = '04:46:27'
$ws.Cells.Item(28, "B").Value = '05:54'
$ws.Cells.Item(28, "C").Value = '10_OLMOS'
$ws.Cells.Item(28, "D").Value = 68
$ws.Cells.Item(28, "E").Value = 'LP1912'

# row 29: 16_SANTA ANA
$ws.Cells.Item(29, "A").Value = '04:46:27'
$ws.Cells.Item(29, "B").Value = '06:04'
$ws.Cells.Item(29, "C").Value = '16_SANTA ANA'
$ws.Cells.Item(29, "D").Value = 78
$ws.Cells.Item(29, "E").Value = 'LP1912'

# row 30: 16_SANTA ANA
$ws.Cells.Item(30, "A").Value = '04:24:09'
$ws.Cells.Item(30, "B").Value = '06:09'
$ws.Cells.Item(30, "C").Value = '16_SANTA ANA'
$ws.Cells.Item(30, "D").Value = 105
$ws.Cells.Item(30, "E").Value = 'LP1912'

# row 31: 215A_EL PATO
$ws.Cells.Item(31, "A").Value = '04:46:27'
$ws.Cells.Item(31, "B").Value = '06:11'
$ws.Cells.Item(31, "C").Value = '215A_EL PATO'
$ws.Cells.Item(31, "D").Value = 85
$ws.Cells.Item(31, "E").Value = 'LP1912'

# row 32: 225_HARAS DEL SUR
$ws.Cells.Item(32, "A").Value = '04:46:27'
$ws.Cells.Item(32, "B").Value = '06:14'
$ws.Cells.Item(32, "C").Value = '225_HARAS DEL SUR'
$ws.Cells.Item(32, "D").Value = 88
$ws.Cells.Item(32, "E").Value = 'LP1912'

# row 33: 26_HERNANDEZ
$ws.Cells.Item(33, "A").Value = '04:46:27'
$ws.Cells.Item(33, "B").Value = '06:21'
$ws.Cells.Item(33, "C").Value = '26_HERNANDEZ'
$ws.Cells.Item(33, "D").Value = 95
$ws.Cells.Item(33, "E").Value = 'LP1912'

# row 34: 23_HERNANDEZ
$ws.Cells.Item(34, "A").Value = '04:46:27'
$ws.Cells.Item(34, "B").Value = '06:27'
$ws.Cells.Item(34, "C").Value = '23_HERNANDEZ'
$ws.Cells.Item(34, "D").Value = 101
$ws.Cells.Item(34, "E").Value = 'LP1912'

# row 35: 86_EST CHICA-ESC AGRARIA
$ws.Cells.Item(35, "A").Value = '04:46:27'
$ws.Cells.Item(35, "B").Value = '06:29'
$ws.Cells.Item(35, "C").Value = '86_EST CHICA-ESC AGRARIA'
$ws.Cells.Item(35, "D").Value = 103
$ws.Cells.Item(35, "E").Value = 'LP1912'

# row 36: 16_SANTA ANA
$ws.Cells.Item(36, "A").Value = '04:46:27'
$ws.Cells.Item(36, "B").Value = '06:31'
$ws.Cells.Item(36, "C").Value = '16_SANTA ANA'
$ws.Cells.Item(36, "D").Value = 105
$ws.Cells.Item(36, "E").Value = 'LP1912'

# row 37: 225_C ROCA-H SUR
$ws.Cells.Item(37, "A").Value = '04:46:27'
$ws.Cells.Item(37, "B").Value = '06:44'
$ws.Cells.Item(37, "C").Value = '225_C ROCA-H SUR'
$ws.Cells.Item(37, "D").Value = 118
$ws.Cells.Item(37, "E").Value = 'LP1912'


# ===== Sheet 2: LP1912-215 =====
$ws = $wb.Worksheets.Item(2)

# Refresh timestamp / row-count banner
$ws.Cells.Item(2, "A").Value = 'Última actualización: 04:46:27'
$ws.Cells.Item(3, "A").Value = 'Total filas: 10'

# Data rows (Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada)
# row 11: 215A_EL PATO
$ws.Cells.Item(11, "A").Value = '04:46:27'
$ws.Cells.Item(11, "B").Value = '04:46'
$ws.Cells.Item(11, "C").Value = '215A_EL PATO'
$ws.Cells.Item(11, "D").Value = 0
$ws.Cells.Item(11, "E").Value = 'LP1912'

# row 12: 215B_EL PATO
$ws.Cells.Item(12, "A").Value = '04:46:27'
$ws.Cells.Item(12, "B").Value = '05:34'
$ws.Cells.Item(12, "C").Value = '215B_EL PATO'
$ws.Cells.Item(12, "D").Value = 48
$ws.Cells.Item(12, "E").Value = 'LP1912'

# row 13: 215B_EL PATO
$ws.Cells.Item(13, "A").Value = '04:24:09'
$ws.Cells.Item(13, "B").Value = '05:35'
$ws.Cells.Item(13, "C").Value = '215B_EL PATO'
$ws.Cells.Item(13, "D").Value = 71
$ws.Cells.Item(13, "E").Value = 'LP1912'

# row 14: 215A_EL PATO
$ws.Cells.Item(14, "A").Value = '04:24:09'
$ws.Cells.Item(14, "B").Value = '05:48'
$ws.Cells.Item(14, "C").Value = '215A_EL PATO'
$ws.Cells.Item(14, "D").Value = 84
$ws.Cells.Item(14, "E").Value = 'LP1912'

# row 15: 215A_EL PATO
$ws.Cells.Item(15, "A").Value = '04:46:27'
$ws.Cells.Item(15, "B").Value = '06:11'
$ws.Cells.Item(15, "C").Value = '215A_EL PATO'
$ws.Cells.Item(15, "D").Value = 85
$ws.Cells.Item(15, "E").Value = 'LP1912'


# ===== Sheet 3: 6203-6173 =====
$ws = $wb.Worksheets.Item(3)

# Refresh timestamp / row-count banner
$ws.Cells.Item(2, "A").Value = 'Última actualización: 04:46:27'
$ws.Cells.Item(3, "A").Value = 'Total filas: 4'

# Data rows (Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada)
# row 7: 215A_LA PLATA
$ws.Cells.Item(7, "A").Value = '04:46:27'
$ws.Cells.Item(7, "B").Value = '05:44'
$ws.Cells.Item(7, "C").Value = '215A_LA PLATA'
$ws.Cells.Item(7, "D").Value = 58
$ws.Cells.Item(7, "E").Value = 'L6173'

# row 8: 215A_LA PLATA
$ws.Cells.Item(8, "A").Value = '04:46:27'
$ws.Cells.Item(8, "B").Value = '06:09'
$ws.Cells.Item(8, "C").Value = '215A_LA PLATA'
$ws.Cells.Item(8, "D").Value = 83
$ws.Cells.Item(8, "E").Value = 'L6173'

# row 9: 215C_LA PLATA
$ws.Cells.Item(9, "A").Value = '04:46:27'
$ws.Cells.Item(9, "B").Value = '06:33'
$ws.Cells.Item(9, "C").Value = '215C_LA PLATA'
$ws.Cells.Item(9, "D").Value = 107
$ws.Cells.Item(9, "E").Value = 'L6203'
